# "1 month results" update:
#  - rename the "Residual" header (H1) to "Error"
#  - add new analysis columns: Relative Error (I, already existed), Squared
#    Error (J), Abs Error (K), plus a small side panel of summary stats
#    (Sum Squared Error / RMSE in M1:N1, Sum Abs Error / RMAE labels in
#    M3:O3, RRMSE in O1) with per-row helper formulas in M2:O2 and M4:O4
#  - fill in row 5's previously-missing Actual value and extend the Error /
#    Relative Error / Squared Error / Abs Error formulas into it
#  - append two new weekly data rows (6 and 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row + summary labels ------------------------------------------
# (written in this order because it is what determines the order new
#  strings land in the shared-strings table)
$ws.Range("H1").Value = "Error"
$ws.Range("I1").Value = "Relative Error"
$ws.Range("J1").Value = "Squared Error"
$ws.Range("M1").Value = "Sum Squared Error"
$ws.Range("N1").Value = "RMSE"
$ws.Range("K1").Value = "Abs Error"
$ws.Range("M3").Value = "Sum Abs Error"
$ws.Range("N3").Value = "RMAE"
$ws.Range("O3").Value = "RMAE"
$ws.Range("O1").Value = "RRMSE"

# --- Row 2: typed individually (not part of a fill-down) ------------------
$ws.Range("J2").Formula = "=H2^2"
$ws.Range("K2").Formula = "=ABS(H2)"
$ws.Range("M2").Formula = "=SUM(J:J)"
$ws.Range("N2").Formula = "=SQRT(M2/(COUNTA(G:G)-1))"
$ws.Range("O2").Formula = "=N2/AVERAGE(G:G)"

# --- Row 5: fill the previously-incomplete row (H5/I5 formulas are added
#     below together with row 6, since they end up as a shared-formula
#     range spanning H5:H6 / I5:I6) ------------------------------------------
$ws.Range("G5").Value = 6.49

# --- Row 6 (new) ------------------------------------------------------------
$ws.Range("A6").Value = 5
$ws.Range("C6").Formula = "=C5+7"
$ws.Range("D6").Value = 5.9211803567514103
$ws.Range("E6").Value = 6.2768796432486296
$ws.Range("F6").Value = 6.09903
$ws.Range("G6").Value = 6.33

# --- Row 7 (new) ------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("C7").Formula = "=C6+7"
$ws.Range("D7").Value = 6.0013503567514102
$ws.Range("E7").Value = 6.3570496432486303
$ws.Range("F7").Value = 6.1791999999999998

# --- Fill-down ranges: these become Excel "shared formulas" because a
#     single relative formula is dropped across a contiguous multi-cell
#     range in one shot (mirrors a select + Ctrl+D fill-down) -------------
$ws.Range("H5:H6").Formula = "=G5-F5"
$ws.Range("I5:I6").Formula = "=H5/G5"
$ws.Range("J3:J6").Formula = "=H3^2"
$ws.Range("K3:K6").Formula = "=ABS(H3)"
$ws.Range("B6:B7").Formula = "=B5+7"

# --- Row 4 summary stats (typed individually) ------------------------------
$ws.Range("M4").Formula = "=SUM(K:K)"
$ws.Range("N4").Formula = "=M4/(COUNTA(G:G)-1)"
$ws.Range("O4").Formula = "=N4/AVERAGE(G:G)"

# --- Apply the m/d date format to the new date cells so they match B2:C5 --
# (copy/paste-formats-only so the existing date style is reused instead of
#  Excel minting a brand-new custom numFmt)
$ws.Range("B5:C5").Copy() | Out-Null
$ws.Range("B6:C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column widths: re-run best-fit only on the columns whose content
#     actually changed width ------------------------------------------------
$ws.Columns("B:E").AutoFit() | Out-Null
$ws.Columns("I:J").AutoFit() | Out-Null
$ws.Columns("M:M").AutoFit() | Out-Null

# --- Selection / view state matching the saved file -------------------------
$ws.Range("O4").Select()
